$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.101.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.47%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.252.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.82%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'607.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.24%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'156.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.04%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  +0.24%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.250.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.64%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.04%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'  -0.19%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'5.67"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.59%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.491"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.54%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.0000265"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.92%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'38.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.60%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'3.789.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.02%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'66.183.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.31%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.255.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.99%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'7.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.17%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  +1.13%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'497.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.51%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'15.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.61%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.744"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.46%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'8.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.75%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'14.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.09%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'86.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.52%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  -0.14%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +0.43%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'9.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.61%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.33%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.131"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +41.93%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'7.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.31%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'2.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -7.75%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'27.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.02%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.02%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.46%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'6.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.30%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'3.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +14.74%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'55.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.29%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'490.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.62%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +4.82%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.0419"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.27%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.129"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.49%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'8.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.76%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'2.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.84%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'2.988.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.29%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.289"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.68%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'28.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.46%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'2.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.38%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  +1.49%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D51").Value = "'121.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.87%  "
$ws.Range("E51").Style = "Normal"
